$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 42608.890243055554
$ws.Range("B6").Value = -4
$ws.Range("C6").Value = 53
$ws.Range("D6").Value = 44
$ws.Range("E6").Value = 35
$ws.Range("F6").Value = 64
$ws.Range("G6").Value = 24724
$ws.Range("H6").Value = 16568
$ws.Range("I6").Value = 2819
$ws.Range("J6").Value = 233
$ws.Range("K6").Value = 192
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 11
$ws.Range("N6").Value = "Noun"
